$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in A3: "gunction" -> "function"
$ws.Range("A3").Value = "Madeline, function measurments"

# Continue the meeting log with the progress made before the meeting / figure 1
$ws.Range("A14").Value = "Madeline, fig 2 regresssions of agg Aug data"
$ws.Range("B14").Value = "trend looks more yearly, rather than a jump in 2020"

$ws.Range("A15").Value = "Yang, Matlab splines"
$ws.Range("B15").Value = 'fit splines in matlab, "not much diff than Madeline"'

$ws.Range("A16").Value = "Nick, R splines"
$ws.Range("B16").Value = "fit splines in r,working on roots, to get function measures."

$ws.Range("A17").Value = "Madeline, differing feature is trend with past trends"

$ws.Range("A18").Value = "Nick, what about weekend, or weekend "

$ws.Range("B19").Value = "This week's theme: Improve on last weeks ideas."
$ws.Range("B19").Font.Bold = $true

# Column B narrows now that it no longer needs to fit the old wide text
$ws.Columns("B").ColumnWidth = 48.65

$ws.Range("B19").Select()
